$d = $word.ActiveDocument

# Useful characters
$rsq = [char]0x2019   # right single quotation mark (curly apostrophe)

# ---------------------------------------------------------------------
# 1) Remove the whole paragraph "Some relevant results were brought to
#    our attention..." (it gets deleted entirely; the following
#    paragraph "We started with..." becomes the merged paragraph).
# ---------------------------------------------------------------------
$target1 = "Some relevant results were brought to our attention and help us to understand what" + $rsq + "s the best way to improve a student" + $rsq + "s scores and help him to live better the scholastic field. "
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Some relevant results were brought to our attention")) {
        $p.Range.Delete()
        $found = $true
        break
    }
}
Write-Host "Removed intro paragraph:" $found

# ---------------------------------------------------------------------
# 2) "we find out that, as we expected, " -> "we find out, as we expected, that "
# ---------------------------------------------------------------------
$r = $d.Content
$ok2 = $r.Find.Execute(
    "we find out that, as we expected, ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "we find out, as we expected, that ", 2)
Write-Host "Step2:" $ok2

# ---------------------------------------------------------------------
# 3) "...influence its students' results, " -> "...influences its students' results, "
# ---------------------------------------------------------------------
$r = $d.Content
$ok3 = $r.Find.Execute(
    "influence its students" + $rsq + " results, ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "influences its students" + $rsq + " results, ", 2)
Write-Host "Step3:" $ok3

# ---------------------------------------------------------------------
# 4) "...evidence in the interaction " -> "...evidence of this in the coefficient of the interaction "
# ---------------------------------------------------------------------
$r = $d.Content
$ok4 = $r.Find.Execute(
    "evidence in the interaction ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "evidence of this in the coefficient of the interaction ", 2)
Write-Host "Step4:" $ok4

# ---------------------------------------------------------------------
# 5) Socioeconomic sentence rewording
# ---------------------------------------------------------------------
$r = $d.Content
$old5 = "given mostly by the socioeconomic state of the student: for both native and immigrant students, we can highlight that the scores in math and reading improves for students with a better socioeconomic situation at home"
$new5 = "given mostly by the socioeconomic state of the student, both for native and immigrant students"
$ok5 = $r.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2)
Write-Host "Step5:" $ok5

# ---------------------------------------------------------------------
# 6) Move <w:lastRenderedPageBreak/> from "Bullied" run. The object
#    model has no direct way to target that internal element, but
#    touching (and then restoring) the run's text forces the engine to
#    rebuild the run without it.
# ---------------------------------------------------------------------
$full = $d.Content.Text
$idxB = $full.IndexOf("Bullied")
if ($idxB -ge 0) {
    $rb = $d.Range($idxB, $idxB + "Bullied".Length)
    $rb.Text = "BullieD"
    $full2 = $d.Content.Text
    $idxB2 = $full2.IndexOf("BullieD")
    $rb2 = $d.Range($idxB2, $idxB2 + "BullieD".Length)
    $rb2.Text = "Bullied"
    Write-Host "Step6: cleared lastRenderedPageBreak on Bullied run"
}
